$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the three "Requisitos" entries so that the
# "LOM3246 - Técnicas de Caracterização de Materiais (Indicação de Conjunto)"
# row moves from the first position to the last position, shifting the
# other two rows up by one.

$reqLOM3246 = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$reqLOB1021 = "LOB1021 -  Física IV  (Requisito)`n"
$reqLOM3016 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

# New order: LOB1021, LOM3016, LOM3246
$ws.Range("B23").Value = $reqLOB1021
$ws.Range("C23").Value = $reqLOB1021

$ws.Range("B24").Value = $reqLOM3016
$ws.Range("C24").Value = $reqLOM3016

$ws.Range("B25").Value = $reqLOM3246
$ws.Range("C25").Value = $reqLOM3246
